$d = $word.ActiveDocument

# 1. Créditos-trabalho: 2 -> 0
$d.Content.Find.Execute("Créditos-trabalho: 2", $true, $false, $false, $false, $false, $true, 1, $false, "Créditos-trabalho: 0", 2)

# 2. Carga horária: 90 h -> 30 h
$d.Content.Find.Execute("Carga horária: 90 h", $true, $false, $false, $false, $false, $true, 1, $false, "Carga horária: 30 h", 2)

# 3. Ativação: 01/01/2024 -> 01/01/2025
$d.Content.Find.Execute("Ativação: 01/01/2024", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2025", 2)

# 4. Add new responsible lecturer before the existing one, same bullet paragraph,
#    separated by a manual line break.
$rng = $d.Content
$rng.Find.Execute("5840535 - Messias Borges Silva", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lineBreak = [char]11
$rng.InsertBefore("3295113 - José Eduardo Holler Branco" + $lineBreak)

# 5. Update evaluation method text
$d.Content.Find.Execute("Provas, trabalhos em grupo, exercícios individuais e seminários.", $true, $false, $false, $false, $false, $true, 1, $false, "Aulas Expositivas; trabalhos e seminários", 2)
